$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "No Flight Can Truly Be COVID-Free"
$ws.Range("B2").Value = "https://lifehacker.com/no-flight-can-truly-be-covid-free-1845898434"
$ws.Range("A3").Value = "We Have Two COVID Vaccines Now"
$ws.Range("B3").Value = "https://vitals.lifehacker.com/we-have-two-covid-vaccines-now-1845913299"
$ws.Range("A4").Value = "What’s Up With the Latest COVID Relief Bill?"
$ws.Range("B4").Value = "https://twocents.lifehacker.com/what-s-up-with-the-latest-covid-relief-bill-1845827491"
$ws.Range("A5").Value = "What We Know About Allergic Reactions to the COVID Vaccines"
$ws.Range("B5").Value = "https://vitals.lifehacker.com/what-we-know-about-allergic-reactions-to-the-covid-vacc-1845934680"
$ws.Range("A6").Value = "What to Do If Youre Exposed to COVID-19"
$ws.Range("B6").Value = "https://lifehacker.com/what-to-do-if-youre-exposed-to-covid-19-1845860079"
$ws.Range("A7").Value = "What Will It Feel Like to Get a COVID Vaccine?"
$ws.Range("B7").Value = "https://vitals.lifehacker.com/what-will-it-feel-like-to-get-a-covid-vaccine-1845810703"
$ws.Range("A8").Value = "Were Probably Days Away From a COVID Vaccine"
$ws.Range("B8").Value = "https://vitals.lifehacker.com/were-probably-days-away-from-a-covid-vaccine-1845861896"
$ws.Range("A9").Value = "Dont Spit on the Ground During a Pandemic (Or Ever)"
$ws.Range("B9").Value = "https://lifehacker.com/dont-spit-on-the-ground-during-a-pandemic-or-ever-1845914513"
$ws.Range("A10").Value = "Californias COVID-19 exposure notification app starts rolling out"
$ws.Range("B10").Value = "https://www.engadget.com/california-covid-19-exposure-notification-app-235001952.html"
$ws.Range("A11").Value = "The Free Market Approach to This Pandemic Isnt Working"
$ws.Range("B11").Value = "https://www.wired.com/story/the-free-market-approach-to-this-pandemic-isnt-working/"
$ws.Range("A12").Value = "Facebook adds new notifications for COVID-19 misinformation"
$ws.Range("B12").Value = "https://www.engadget.com/facebook-notifications-harmful-coronavirus-misinformation-191233820.html"
$ws.Range("A13").Value = "Congress approves COVID-19 spending bill with contentious copyright measures"
$ws.Range("B13").Value = "https://www.engadget.com/covid-19-spending-bill-passes-with-new-streaming-copyright-law-tacked-on-102046838.html"
$ws.Range("A14").Value = "Twitter will start removing COVID-19 vaccine misinformation next week"
$ws.Range("B14").Value = "https://www.engadget.com/twitter-covid-19-vaccine-misinformation-policy-204452958.html"
$ws.Range("A15").Value = "Covid-19: How Covid cruise ships are navigating troubled waters"
$ws.Range("B15").Value = "https://www.bbc.co.uk/news/av/world-55241333"
$ws.Range("A16").Value = "The Last, ‘Ultra-Cold’ Mile for Covid-19 Vaccines"
$ws.Range("B16").Value = "https://www.wired.com/story/the-last-ultra-cold-mile-for-covid-19-vaccines/"
$ws.Range("A17").Value = "Conferences After Covid Will Be Shorter—and Smarter"
$ws.Range("B17").Value = "https://www.wired.com/story/what-conferences-will-look-like-post-covid/"
$ws.Range("A18").Value = "The ‘Healthy Building’ Surge Will Outlast the Pandemic"
$ws.Range("B18").Value = "https://www.wired.com/story/healthy-building-outlast-pandemic/"
$ws.Range("A19").Value = "A Clever Strategy to Distribute Covid Aid—With Satellite Data"
$ws.Range("B19").Value = "https://www.wired.com/story/clever-strategy-distribute-covid-aid-satellite-data/"
$ws.Range("A20").Value = "Poland’s GeneMe secures €5.2M seed funding for its rapid COVID-19 test"
$ws.Range("B20").Value = "http://techcrunch.com/2020/12/17/polands-geneme-secures-e5-2m-seed-funding-for-its-rapid-covid-19-test/"
$ws.Range("A21").Value = "Twitter says it will start removing COVID-19 vaccine misinformation"
$ws.Range("B21").Value = "https://www.theverge.com/2020/12/16/22179074/twitter-coronavirus-misinformation-covid19-vaccine-vaccination-label"
